$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style (bold, border, centered) from F1 to the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Update existing numeric values in row 2
$ws.Range("B2").Value = 0.2477729233431636
$ws.Range("C2").Value = 0.9951572719822364
$ws.Range("D2").Value = 0.3903416901953408

# Update the model description text (now wraps to a new line, n_estimators changed 50 -> 150)
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=3, n_estimators=150))])"

# New value cells for the added columns
$ws.Range("G2").Value = 0.1239050709499376
$ws.Range("H2").Value = 0.991
